$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was row data, now row 4 data per permutation)
$ws.Range("D2").Value = 44491
$ws.Range("J2").Value = 150

# Row 3 (was row data, now row 20 data per permutation)
$ws.Range("D3").Value = 44497
$ws.Range("J3").Value = 250

# Row 4 (was row data, now row 9 data per permutation)
$ws.Range("D4").Value = 44312
$ws.Range("J4").Value = 180
$ws.Range("K4").Value = 2500
$ws.Range("L4").Value = 2500
$ws.Range("M4").Value = 2500
$ws.Range("N4").Value = '$/unidad'
$ws.Range("P4").Value = 2500

# Row 5 (was row data, now row 16 data per permutation)
$ws.Range("D5").Value = 44223
$ws.Range("H5").Value = 'Americana O Klondike'
$ws.Range("I5").Value = 'Extra'
$ws.Range("J5").Value = 340
$ws.Range("K5").Value = 2500
$ws.Range("L5").Value = 2500
$ws.Range("M5").Value = 2500
$ws.Range("P5").Value = 2500

# Row 6 (was row data, now row 17 data per permutation)
$ws.Range("D6").Value = 44223
$ws.Range("H6").Value = 'Americana O Klondike'
$ws.Range("I6").Value = 'Primera'
$ws.Range("J6").Value = 400
$ws.Range("K6").Value = 2000
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 2000
$ws.Range("P6").Value = 2000

# Row 7 (was row data, now row 18 data per permutation)
$ws.Range("D7").Value = 44223
$ws.Range("H7").Value = 'Americana O Klondike'
$ws.Range("I7").Value = 'Segunda'
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = 1500
$ws.Range("P7").Value = 1500

# Row 8 (was row data, now row 19 data per permutation)
$ws.Range("D8").Value = 44223
$ws.Range("H8").Value = 'Americana O Klondike'
$ws.Range("I8").Value = 'Tercera'
$ws.Range("J8").Value = 160
$ws.Range("K8").Value = 1000
$ws.Range("L8").Value = 1000
$ws.Range("M8").Value = 1000
$ws.Range("O8").Value = 'Región de O''Higgins'
$ws.Range("P8").Value = 1000

# Row 9 (was row data, now row 5 data per permutation)
$ws.Range("D9").Value = 44167
$ws.Range("J9").Value = 400
$ws.Range("K9").Value = 5000
$ws.Range("L9").Value = 5000
$ws.Range("M9").Value = 5000
$ws.Range("O9").Value = 'Región de O''Higgins'
$ws.Range("P9").Value = 5000

# Row 10 (was row data, now row 6 data per permutation)
$ws.Range("D10").Value = 44167
$ws.Range("I10").Value = 'Segunda'
$ws.Range("J10").Value = 560
$ws.Range("K10").Value = 3000
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 3000
$ws.Range("N10").Value = '$/unidad'
$ws.Range("O10").Value = 'Región de O''Higgins'
$ws.Range("P10").Value = 3000

# Row 11 (was row data, now row 7 data per permutation)
$ws.Range("D11").Value = 44167
$ws.Range("I11").Value = 'Tercera'
$ws.Range("J11").Value = 450
$ws.Range("K11").Value = 2000
$ws.Range("L11").Value = 2000
$ws.Range("M11").Value = 2000
$ws.Range("N11").Value = '$/unidad'
$ws.Range("O11").Value = 'Región de O''Higgins'
$ws.Range("P11").Value = 2000

# Row 12 (was row data, now row 10 data per permutation)
$ws.Range("D12").Value = 44477
$ws.Range("I12").Value = 'Primera'
$ws.Range("J12").Value = 80
$ws.Range("K12").Value = 800
$ws.Range("L12").Value = 800
$ws.Range("M12").Value = 800
$ws.Range("N12").Value = '$/kilo (volumen en unidades)'
$ws.Range("O12").Value = 'Perú'
$ws.Range("P12").Value = 800

# Row 13 (was row data, now row 11 data per permutation)
$ws.Range("D13").Value = 44488
$ws.Range("J13").Value = 150
$ws.Range("K13").Value = 800
$ws.Range("L13").Value = 800
$ws.Range("M13").Value = 800
$ws.Range("N13").Value = '$/kilo (volumen en unidades)'
$ws.Range("O13").Value = 'Perú'
$ws.Range("P13").Value = 800

# Row 14 (was row data, now row 8 data per permutation)
$ws.Range("D14").Value = 44305
$ws.Range("I14").Value = 'Primera'
$ws.Range("J14").Value = 100
$ws.Range("O14").Value = 'Perú'

# Row 15 (was row data, now row 14 data per permutation)
$ws.Range("I15").Value = 'Extra'
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 2500
$ws.Range("L15").Value = 2500
$ws.Range("M15").Value = 2500
$ws.Range("P15").Value = 2500

# Row 16 (was row data, now row 15 data per permutation)
$ws.Range("D16").Value = 44217
$ws.Range("H16").Value = 'Sin especificar'
$ws.Range("I16").Value = 'Primera'
$ws.Range("J16").Value = 280
$ws.Range("K16").Value = 2000
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = 2000
$ws.Range("P16").Value = 2000

# Row 17 (was row data, now row 21 data per permutation)
$ws.Range("D17").Value = 44504
$ws.Range("H17").Value = 'Sin especificar'
$ws.Range("J17").Value = 200
$ws.Range("K17").Value = 800
$ws.Range("L17").Value = 800
$ws.Range("M17").Value = 800
$ws.Range("N17").Value = '$/kilo (volumen en unidades)'
$ws.Range("O17").Value = 'Perú'
$ws.Range("P17").Value = 800

# Row 18 (was row data, now row 2 data per permutation)
$ws.Range("D18").Value = 44483
$ws.Range("H18").Value = 'Sin especificar'
$ws.Range("I18").Value = 'Primera'
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = 800
$ws.Range("L18").Value = 800
$ws.Range("M18").Value = 800
$ws.Range("N18").Value = '$/kilo (volumen en unidades)'
$ws.Range("O18").Value = 'Perú'
$ws.Range("P18").Value = 800

# Row 19 (was row data, now row 3 data per permutation)
$ws.Range("D19").Value = 44495
$ws.Range("H19").Value = 'Sin especificar'
$ws.Range("I19").Value = 'Primera'
$ws.Range("J19").Value = 200
$ws.Range("K19").Value = 800
$ws.Range("L19").Value = 800
$ws.Range("M19").Value = 800
$ws.Range("N19").Value = '$/kilo (volumen en unidades)'
$ws.Range("O19").Value = 'Perú'
$ws.Range("P19").Value = 800

# Row 20 (was row data, now row 22 data per permutation)
$ws.Range("D20").Value = 44510

# Row 21 (was row data, now row 12 data per permutation)
$ws.Range("D21").Value = 44194
$ws.Range("I21").Value = 'Extra'
$ws.Range("J21").Value = 120
$ws.Range("K21").Value = 3500
$ws.Range("L21").Value = 3500
$ws.Range("M21").Value = 3500
$ws.Range("N21").Value = '$/unidad'
$ws.Range("O21").Value = 'Región de O''Higgins'
$ws.Range("P21").Value = 3500

# Row 22 (was row data, now row 13 data per permutation)
$ws.Range("D22").Value = 44194
$ws.Range("J22").Value = 200
$ws.Range("K22").Value = 3000
$ws.Range("L22").Value = 3000
$ws.Range("M22").Value = 3000
$ws.Range("N22").Value = '$/unidad'
$ws.Range("O22").Value = 'Región de O''Higgins'
$ws.Range("P22").Value = 3000
